$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Invalid block: Unexpected tag EOF missing [ENDFOR] while parsing m:for v | self.eClassifiers",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "    <---Invalid block: Unexpected tag EOF missing [ENDFOR] while parsing m:for v | self.eClassifiers",
    2
)
